$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill previously-blank score cells (existing style retained) ---
$ws.Range("C3:I3").Value = 1
$ws.Range("G6:I6").Value = 2
$ws.Range("H8:I8").Value = 2
$ws.Range("G9:I9").Value = 2
$ws.Range("H13:I13").Value = 2
$ws.Range("G17:I17").Value = 2
$ws.Range("G24:I24").Value = 2
$ws.Range("G25:I25").Value = 2
$ws.Range("H26").Value = 2
$ws.Range("G27:I27").Value = 2
$ws.Range("I29").Value = 2
$ws.Range("H30:I30").Value = 2
$ws.Range("H32:I32").Value = 2

# --- New J-column cells: set value, then copy formatting (border/fill) from a sibling cell that already has the desired style ---
$ws.Range("J3").Value = 1
$ws.Range("J2").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("J6").Value = 2
$ws.Range("J4").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("J7").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Range("J8").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("J9").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$ws.Range("J10").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J10").PasteSpecial(-4122)
$ws.Range("J11").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J11").PasteSpecial(-4122)
$ws.Range("J13").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J13").PasteSpecial(-4122)
$ws.Range("J17").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J17").PasteSpecial(-4122)
$ws.Range("J21").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J21").PasteSpecial(-4122)
$ws.Range("J24").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J24").PasteSpecial(-4122)
$ws.Range("J25").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J25").PasteSpecial(-4122)
$ws.Range("J27").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J27").PasteSpecial(-4122)
$ws.Range("J28").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J28").PasteSpecial(-4122)
$ws.Range("J29").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J29").PasteSpecial(-4122)
$ws.Range("J30").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J30").PasteSpecial(-4122)
$ws.Range("J31").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J31").PasteSpecial(-4122)
$ws.Range("J32").Value = 2
$ws.Range("J18").Copy()
$ws.Range("J32").PasteSpecial(-4122)

# --- Update active selection to J3 (matches commit's recorded cursor position) ---
$null = $ws.Range("J3").Select()
